$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.810.56'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.707.00'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '678.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.77%  '
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E11').Value = '  +2.20%  '
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '33.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').Value = '3.718.47'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').Value = '69.780.19'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '473.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '80.50'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').Value = '3.854.49'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E24').Value = '  +3.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.65'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('E34').Value = '  +4.14%  '
$ws.Range('D35').Value = '3.696.24'
$ws.Range('E35').Value = '  +0.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '8.61'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.03%  '
$ws.Range('E37').Value = '  +1.23%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0910'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '168.06'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.48%  '
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('E47').Value = '  +2.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.269'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.81%  '
